$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Create the even-page and first-page header/footer variants -------------
# (Word splits a single default header/footer into even/default/first
#  variants once you touch the even- or first-page story; this mirrors the
#  headerReference/footerReference w:type="even"/"first" entries added to
#  sectPr in the target document.)
$hEven = $sec.Headers(3)
$hEven.Range.Text = ""
$hFirst = $sec.Headers(2)
$hFirst.Range.Text = ""

$fEven = $sec.Footers(3)
$fEven.Range.Text = ""
$fFirst = $sec.Footers(2)
$fFirst.Range.Text = ""

# --- Update the FOI office contact details in the default footer ------------
$fDefault = $sec.Footers(1)

$r1 = $fDefault.Range.Find.Execute("V8W 9K1", $false, $false, $false, $false, $false, $true, 1, $false, "V8W 9V1", 2)
Write-Host "postal code updated:" $r1

$r2 = $fDefault.Range.Find.Execute("250-387-9843", $false, $false, $false, $false, $false, $true, 1, $false, "250-387-1321", 2)
Write-Host "phone number updated:" $r2

# --- Drop the stale pagination cache on the "Payment Method:" run -----------
$r3 = $d.Content.Find.Execute("Payment Method: ", $false, $false, $false, $false, $false, $true, 1, $false, "Payment Method: ", 2)
Write-Host "payment method run refreshed:" $r3
